$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D), Volume 1h (E), and Hora (G) columns for rows 2-51.
# Values are written with a leading apostrophe so Excel stores them as text
# (matching the inlineStr cells already in the workbook) and the cell Style is
# reset to "Normal" afterward so no new number-format style gets introduced.

$ws.Range("D2").Value = "'304.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'4.13%"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "'19"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'35.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'14.64%"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Value = "'19"
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = "'5.087"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'2.43%"
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").Value = "'19"
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").Value = "'0.07833"
$ws.Range("D5").Style = "Normal"
$ws.Range("G5").Value = "'19"
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").Value = "'2.267"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.76%"
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").Value = "'19"
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = "'8.121"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'4.37%"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").Value = "'19"
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = "'4.006"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'6.09%"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "'19"
$ws.Range("G8").Style = "Normal"
$ws.Range("D9").Value = "'0.9258"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.64%"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").Value = "'19"
$ws.Range("G9").Style = "Normal"
$ws.Range("D10").Value = "'0.09864"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'6.09%"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Value = "'19"
$ws.Range("G10").Style = "Normal"
$ws.Range("D11").Value = "'0.1825"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'5.25%"
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").Value = "'19"
$ws.Range("G11").Style = "Normal"
$ws.Range("D12").Value = "'0.08710"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'4.28%"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").Value = "'19"
$ws.Range("G12").Style = "Normal"
$ws.Range("D13").Value = "'0.03397"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'3.29%"
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").Value = "'19"
$ws.Range("G13").Style = "Normal"
$ws.Range("D14").Value = "'0.09934"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.09%"
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").Value = "'19"
$ws.Range("G14").Style = "Normal"
$ws.Range("D15").Value = "'0.001475"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.73%"
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").Value = "'19"
$ws.Range("G15").Style = "Normal"
$ws.Range("D16").Value = "'0.005788"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.18%"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").Value = "'19"
$ws.Range("G16").Style = "Normal"
$ws.Range("D17").Value = "'3.480"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.23%"
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").Value = "'19"
$ws.Range("G17").Style = "Normal"
$ws.Range("E18").Value = "'-1.24%"
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").Value = "'19"
$ws.Range("G18").Style = "Normal"
$ws.Range("D19").Value = "'0.3433"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'3.00%"
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").Value = "'19"
$ws.Range("G19").Style = "Normal"
$ws.Range("E20").Value = "'0.51%"
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").Value = "'19"
$ws.Range("G20").Style = "Normal"
$ws.Range("D21").Value = "'4.538"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'10.66%"
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").Value = "'19"
$ws.Range("G21").Style = "Normal"
$ws.Range("D22").Value = "'0.2235"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'5.26%"
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").Value = "'19"
$ws.Range("G22").Style = "Normal"
$ws.Range("D23").Value = "'0.04675"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'3.03%"
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").Value = "'19"
$ws.Range("G23").Style = "Normal"
$ws.Range("E24").Value = "'1.59%"
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").Value = "'19"
$ws.Range("G24").Style = "Normal"
$ws.Range("D25").Value = "'0.004509"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'4.75%"
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").Value = "'19"
$ws.Range("G25").Style = "Normal"
$ws.Range("E26").Value = "'-0.01%"
$ws.Range("E26").Style = "Normal"
$ws.Range("G26").Value = "'19"
$ws.Range("G26").Style = "Normal"
$ws.Range("D27").Value = "'0.0002699"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-20.46%"
$ws.Range("E27").Style = "Normal"
$ws.Range("G27").Value = "'19"
$ws.Range("G27").Style = "Normal"
$ws.Range("G28").Value = "'19"
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").Value = "'19"
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").Value = "'19"
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = "'19"
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = "'19"
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").Value = "'19"
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").Value = "'19"
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").Value = "'19"
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").Value = "'19"
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").Value = "'19"
$ws.Range("G37").Style = "Normal"
$ws.Range("G38").Value = "'19"
$ws.Range("G38").Style = "Normal"
$ws.Range("D39").Value = "'0.01760"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'8.67%"
$ws.Range("E39").Style = "Normal"
$ws.Range("G39").Value = "'19"
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = "'0.04719"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.04%"
$ws.Range("E40").Style = "Normal"
$ws.Range("G40").Value = "'19"
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").Value = "'0.007863"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'4.71%"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").Value = "'19"
$ws.Range("G41").Style = "Normal"
$ws.Range("D42").Value = "'0.1421"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'4.39%"
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").Value = "'19"
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").Value = "'0.008647"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-12.10%"
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").Value = "'19"
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").Value = "'0.002288"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'3.14%"
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").Value = "'19"
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = "'0.009120"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-6.48%"
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").Value = "'19"
$ws.Range("G45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006170"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'1.07%"
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").Value = "'19"
$ws.Range("G46").Style = "Normal"
$ws.Range("E47").Value = "'0.10%"
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").Value = "'19"
$ws.Range("G47").Style = "Normal"
$ws.Range("D48").Value = "'5.675"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'103.07%"
$ws.Range("E48").Style = "Normal"
$ws.Range("G48").Value = "'19"
$ws.Range("G48").Style = "Normal"
$ws.Range("D49").Value = "'0.002690"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'34.55%"
$ws.Range("E49").Style = "Normal"
$ws.Range("G49").Value = "'19"
$ws.Range("G49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.10%"
$ws.Range("E50").Style = "Normal"
$ws.Range("G50").Value = "'19"
$ws.Range("G50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.10%"
$ws.Range("E51").Style = "Normal"
$ws.Range("G51").Value = "'19"
$ws.Range("G51").Style = "Normal"
